$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# Update the "Marking" row correct-answer marks value (B11: 3 -> 5)
$ws.Range("B11").Value = 5

# Update the "Total" row correct marks (B12: 63 -> 105)
$ws.Range("B12").Value = 105

# Update the corr/total marks display string (E12: "62/84" -> "105/140")
$ws.Range("E12").Value = "105/140"
